$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.059.76"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.789.21"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.26"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.296"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0686"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0940"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "2.046.10"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.37"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").Value = "1.787.12"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.623"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "34.048.59"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.89"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.10"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.74"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.10"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.94"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.15"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.20"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +2.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0517"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.65"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.63"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.09%  "
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("D35").Value = "1.401.56"
$ws.Range("E35").Value = "  +1.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.654"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("E39").Value = "  +7.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "79.96"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.919"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.68"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +14.30%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.68"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.12"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.97%  "
$ws.Range("E46").Value = "  +5.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0508"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.57%  "
$ws.Range("E48").Value = "  +2.29%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").Value = "1.948.10"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("E51").Value = "  -0.04%  "
